$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price/volume snapshot cells to match the refreshed feed.
# "Price" (column D) values are stored as text (e.g. thousands separated by
# "." and precise trailing zeros), so NumberFormat is forced to "@" before
# assignment to stop Excel from silently re-typing them as numbers.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.014.12'
$ws.Range('E2').Value = '  +2.21%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.845.93'
$ws.Range('E3').Value = '  +2.00%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '310.28'
$ws.Range('E5').Value = '  +1.19%  '
$ws.Range('E6').Value = '  -0.02%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4664'
$ws.Range('E7').Value = '  +3.40%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07170'
$ws.Range('E9').Value = '  +1.28%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9160'
$ws.Range('E10').Value = '  +2.35%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '19.61'
$ws.Range('E11').Value = '  +1.09%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07694'
$ws.Range('E12').Value = '  -1.50%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.900.24'
$ws.Range('E13').Value = '  +5.00%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.293'
$ws.Range('E14').Value = '  +0.11%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.419'
$ws.Range('E15').Value = '  +1.80%  '
$ws.Range('E16').Value = '  +3.85%  '
$ws.Range('E17').Value = '  +0.03%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008595'
$ws.Range('E18').Value = '  +1.02%  '
$ws.Range('E19').Value = '  +0.05%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '27.043.37'
$ws.Range('E20').Value = '  +2.13%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.40'
$ws.Range('E21').Value = '  +1.51%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.031'
$ws.Range('E22').Value = '  +1.25%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.65'
$ws.Range('E23').Value = '  +1.17%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.935'
$ws.Range('E24').Value = '  -1.24%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '152.60'
$ws.Range('E25').Value = '  +0.49%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '18.42'
$ws.Range('E26').Value = '  +3.49%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.062'
$ws.Range('E27').Value = '  +0.24%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '114.15'
$ws.Range('E28').Value = '  +1.42%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.923'
$ws.Range('E29').Value = '  +1.51%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.08867'
$ws.Range('E30').Value = '  +1.99%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.180'
$ws.Range('E31').Value = '  +2.59%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.857'
$ws.Range('E32').Value = '  +1.89%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.176'
$ws.Range('E33').Value = '  +5.94%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7478'
$ws.Range('E34').Value = '  +2.88%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.473'
$ws.Range('E35').Value = '  +0.49%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.085'
$ws.Range('E36').Value = '  +0.77%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01942'
$ws.Range('E37').Value = '  +0.70%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.965'
$ws.Range('E38').Value = '  +2.75%  '
$ws.Range('E39').Value = '  +0.80%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.5189'
$ws.Range('E40').Value = '  +2.15%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.910'
$ws.Range('E41').Value = '  +2.04%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1510'
$ws.Range('E42').Value = '  -0.49%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.166'
$ws.Range('E43').Value = '  +1.54%  '
$ws.Range('E44').Value = '  +4.55%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.4707'
$ws.Range('E45').Value = '  +0.99%  '
$ws.Range('E46').Value = '  +0.05%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '100.71'
$ws.Range('E47').Value = '  +0.68%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.608'
$ws.Range('E48').Value = '  +2.15%  '
# Row 49/50: coin order swapped (Aave now ranks above Cronos).
$ws.Range('B49').Value = 'Aave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '64.78'
$ws.Range('E49').Value = '  +1.88%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06045'
$ws.Range('E50').Value = '  +0.99%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '36.15'
$ws.Range('E51').Value = '  +0.64%  '
